# This script applies an incremental data refresh to the "广州-漫展信息"
# workbook:
#   - On sheet "展览" the oldest exhibition entry (row 2) was removed; all
#     following rows shift up by one and the last row is dropped.
#   - On sheet "全部类型" the same exhibition entry (row 3 on that sheet,
#     since row 2 holds an unrelated "本地生活" style item) was removed,
#     with the same upward shift and last-row drop.
#   - A handful of "想去人数" (interested-count) figures were bumped on
#     sheets "演出", "本地生活" and "全部类型".
#
# Column A on every sheet is a static, pre-computed row index (0,1,2,...)
# and must stay untouched; only columns B:I carry the data that shifts.

$wb = $excel.ActiveWorkbook

function Remove-FirstDataRow {
    # Positional params: $ws = worksheet ComObject,
    # $deleteRow = row number (1-based, incl. header) whose B:I content is
    # dropped, $lastRow = current last row with data.
    param($ws, $deleteRow, $lastRow)

    # Columns B (plain date, e.g. "2024.02.16") and E (date range, e.g.
    # "2024.02.16 09:30-02.16 16:30") look like dates/times to Excel, so
    # force them to text first to stop COM from silently converting the
    # shifted strings into date serials.
    $bRange = $ws.Range("B$deleteRow`:B$lastRow")
    $eRange = $ws.Range("E$deleteRow`:E$lastRow")
    $bRange.NumberFormat = "@"
    $eRange.NumberFormat = "@"

    # Shift columns B:I up by one row: row N gets what used to be in row N+1.
    $srcRange = $ws.Range("B" + ($deleteRow + 1) + ":I$lastRow")
    $dstRange = $ws.Range("B$deleteRow`:I" + ($lastRow - 1))
    $dstRange.Value2 = $srcRange.Value2

    # Restore the original (General) formatting on B/E now that the text
    # has safely landed, so we don't leave stray number formats behind.
    $bRange.ClearFormats()
    $eRange.ClearFormats()

    # The last row is now a duplicate of row (lastRow - 1); remove it
    # entirely so the sheet's used range/dimension shrinks by one row,
    # just like it does after the source entry was removed upstream.
    $ws.Rows($lastRow).Delete()
}

# --- Sheet "展览": drop the "2024.02.14 广州·运动番only4.0" entry (row 2) ---
$wsExhibit = $wb.Worksheets.Item("展览")
Remove-FirstDataRow $wsExhibit 2 32

# --- Sheet "全部类型": drop the same entry, which sits in row 3 here ---
$wsAll = $wb.Worksheets.Item("全部类型")
Remove-FirstDataRow $wsAll 3 45

# --- Sheet "演出": bump "想去人数" for 春卷饭十周年 (row 10) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F10").Value = 623

# --- Sheet "本地生活": bump "想去人数" for 次元波板糖 (row 2) ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 331

# --- Sheet "全部类型": same 次元波板糖 entry also lives here, untouched
#     by the row shift above since it's row 2 ---
$wsAll.Range("F2").Value = 331
